$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.570987914975575
$ws.Range("C2").Value = 0.1373413713752711
$ws.Range("D2").Value = 0.0763965986979116
$ws.Range("E2").Value = 0.03427439856695003
$ws.Range("G2").Value = 0.002603033726467485
$ws.Range("I2").Value = 2.368000303250767
$ws.Range("K2").Value = 1.499801105360206
$ws.Range("L2").Value = 0.2836558607769746
$ws.Range("N2").Value = 3.77726788383147
$ws.Range("B3").Value = 1.523222088194643
$ws.Range("C3").Value = 0.1264070970914304
$ws.Range("D3").Value = 0.06961488046250963
$ws.Range("E3").Value = 0.03412303381348636
$ws.Range("G3").Value = 0.002608093193199091
$ws.Range("I3").Value = 2.348524274835654
$ws.Range("K3").Value = 1.442416361311984
$ws.Range("L3").Value = 0.2759165171138847
$ws.Range("N3").Value = 3.765180410144566
$ws.Range("B4").Value = 1.49489416250006
$ws.Range("C4").Value = 0.1197698222514134
$ws.Range("D4").Value = 0.06549201376698477
$ws.Range("E4").Value = 0.0340453355302941
$ws.Range("G4").Value = 0.002611362910980765
$ws.Range("I4").Value = 2.337387923276978
$ws.Range("K4").Value = 1.408146798667019
$ws.Range("L4").Value = 0.2713318479904956
$ws.Range("N4").Value = 3.75846783310709
$ws.Range("B5").Value = 1.483601334574672
$ws.Range("C5").Value = 0.1170840324575693
$ws.Range("D5").Value = 0.0638221287052545
$ws.Range("E5").Value = 0.03401750995338126
$ws.Range("G5").Value = 0.002612736520704419
$ws.Range("I5").Value = 2.333056068640204
$ws.Range("K5").Value = 1.394423444462205
$ws.Range("L5").Value = 0.2695055226053
$ws.Range("N5").Value = 3.755910246589224
$ws.Range("B6").Value = 1.481741315149463
$ws.Range("C6").Value = 0.1166391962777311
$ws.Range("D6").Value = 0.06354545892706653
$ws.Range("E6").Value = 0.03401312143354573
$ws.Range("G6").Value = 0.002612967098613671
$ws.Range("I6").Value = 2.332349215876079
$ws.Range("K6").Value = 1.392159269522836
$ws.Range("L6").Value = 0.269204794083052
$ws.Range("N6").Value = 3.755496288479776
$ws.Range("B7").Value = 1.49474084777637
$ws.Range("C7").Value = 0.1197335243509769
$ws.Range("D7").Value = 0.0654694519121648
$ws.Range("E7").Value = 0.03404494472197683
$ws.Range("G7").Value = 0.002611381268973497
$ws.Range("I7").Value = 2.337328667620312
$ws.Range("K7").Value = 1.407960742923962
$ws.Range("L7").Value = 0.2713070477557977
$ws.Range("N7").Value = 3.758432621197073
$ws.Range("B8").Value = 1.554310253631741
$ws.Range("C8").Value = 0.1335552125700019
$ws.Range("D8").Value = 0.07404964546738313
$ws.Range("E8").Value = 0.03421904759935224
$ws.Range("G8").Value = 0.002604744445869289
$ws.Range("I8").Value = 2.361114106790353
$ws.Range("K8").Value = 1.479814035281549
$ws.Range("L8").Value = 0.2809525384716949
$ws.Range("N8").Value = 3.772952631395228
$ws.Range("B9").Value = 1.679097708792312
$ws.Range("C9").Value = 0.1612791036452563
$ws.Range("D9").Value = 0.0912084376433171
$ws.Range("E9").Value = 0.03468123819504321
$ws.Range("G9").Value = 0.002593018079784011
$ws.Range("I9").Value = 2.414302316194593
$ws.Range("K9").Value = 1.628422918963281
$ws.Range("L9").Value = 0.3012012764429386
$ws.Range("N9").Value = 3.807078529464491
$ws.Range("B10").Value = 1.775698442909402
$ws.Range("C10").Value = 0.1820461542558576
$ws.Range("D10").Value = 0.1040286019230621
$ws.Range("E10").Value = 0.0350943284868297
$ws.Range("G10").Value = 0.002585179243632112
$ws.Range("I10").Value = 2.457407964089995
$ws.Range("K10").Value = 1.742384632299206
$ws.Range("L10").Value = 0.3169021094258397
$ws.Range("N10").Value = 3.835636104396571
$ws.Range("B11").Value = 1.820726683132875
$ws.Range("C11").Value = 0.1915847202267855
$ws.Range("D11").Value = 0.1099095811635067
$ws.Range("E11").Value = 0.03529820318210142
$ws.Range("G11").Value = 0.002581779862081879
$ws.Range("I11").Value = 2.477901354746578
$ws.Range("K11").Value = 1.795285199514296
$ws.Range("L11").Value = 0.3242262759416974
$ws.Range("N11").Value = 3.849393405968215
$ws.Range("B12").Value = 1.837934493259979
$ws.Range("C12").Value = 0.1952102218314735
$ws.Range("D12").Value = 0.1121437704218664
$ws.Range("E12").Value = 0.03537769725773288
$ws.Range("G12").Value = 0.002580516408040404
$ws.Range("I12").Value = 2.485789449519629
$ws.Range("K12").Value = 1.81547078995726
$ws.Range("L12").Value = 0.3270260532519984
$ws.Range("N12").Value = 3.854713781584337
$ws.Range("B13").Value = 1.83422150754626
$ws.Range("C13").Value = 0.1944288021523164
$ws.Range("D13").Value = 0.1116622757520247
$ws.Range("E13").Value = 0.03536047494185901
$ws.Range("G13").Value = 0.002580787458279587
$ws.Range("I13").Value = 2.484084919694169
$ws.Range("K13").Value = 1.811116627344461
$ws.Range("L13").Value = 0.3264219003736173
$ws.Range("N13").Value = 3.853563006875817
$ws.Range("B14").Value = 1.822139237577971
$ws.Range("C14").Value = 0.1918827211096072
$ws.Range("D14").Value = 0.1100932443474534
$ws.Range("E14").Value = 0.03530469730936048
$ws.Range("G14").Value = 0.002581675440392985
$ws.Range("I14").Value = 2.478547750915993
$ws.Range("K14").Value = 1.796942801929674
$ws.Range("L14").Value = 0.3244560879438581
$ws.Range("N14").Value = 3.849828892794648
$ws.Range("B15").Value = 1.814758912563548
$ws.Range("C15").Value = 0.1903249337452735
$ws.Range("D15").Value = 0.1091331078698943
$ws.Range("E15").Value = 0.03527083017143084
$ws.Range("G15").Value = 0.002582222452926986
$ws.Range("I15").Value = 2.475172719954202
$ws.Range("K15").Value = 1.788280919153635
$ws.Range("L15").Value = 0.3232553968264398
$ws.Range("N15").Value = 3.847556085066202
$ws.Range("B16").Value = 1.772777615670009
$ws.Range("C16").Value = 0.1814246532145489
$ws.Range("D16").Value = 0.1036452645415977
$ws.Range("E16").Value = 0.03508132560064681
$ws.Range("G16").Value = 0.002585404741217931
$ws.Range("I16").Value = 2.456086513246774
$ws.Range("K16").Value = 1.738948850040344
$ws.Range("L16").Value = 0.3164271272315204
$ws.Range("N16").Value = 3.83475250276652
$ws.Range("B17").Value = 1.747301603238782
$ws.Range("C17").Value = 0.1759882689140113
$ws.Range("D17").Value = 0.1002913057590433
$ws.Range("E17").Value = 0.03496915528938693
$ws.Range("G17").Value = 0.002587399534941125
$ws.Range("I17").Value = 2.444604626204793
$ws.Range("K17").Value = 1.708957142142594
$ws.Range("L17").Value = 0.3122848579057518
$ws.Range("N17").Value = 3.82709460807115
$ws.Range("B18").Value = 1.732750470378335
$ws.Range("C18").Value = 0.1728700020622966
$ws.Range("D18").Value = 0.09836680305932077
$ws.Range("E18").Value = 0.03490614012584281
$ws.Range("G18").Value = 0.002588562571003683
$ws.Range("I18").Value = 2.438083762371008
$ws.Range("K18").Value = 1.691806302035587
$ws.Range("L18").Value = 0.3099194357904622
$ws.Range("N18").Value = 3.822762083330616
$ws.Range("B19").Value = 1.727841203557887
$ws.Range("C19").Value = 0.1718156783070697
$ws.Range("D19").Value = 0.09771598628464062
$ws.Range("E19").Value = 0.03488506239521527
$ws.Range("G19").Value = 0.002588959052649398
$ws.Range("I19").Value = 2.435890186478233
$ws.Range("K19").Value = 1.686016396885805
$ws.Range("L19").Value = 0.3091214776698479
$ws.Range("N19").Value = 3.821307529253772
$ws.Range("B20").Value = 1.750003007072905
$ws.Range("C20").Value = 0.1765660893777294
$ws.Range("D20").Value = 0.1006478627091525
$ws.Range("E20").Value = 0.03498094056118717
$ws.Range("G20").Value = 0.002587185563070547
$ws.Range("I20").Value = 2.445818276704827
$ws.Range("K20").Value = 1.712139495295844
$ws.Range("L20").Value = 0.3127240386236139
$ws.Range("N20").Value = 3.827902338431244
$ws.Range("B21").Value = 1.825683837560348
$ws.Range("C21").Value = 0.1926301993567279
$ws.Range("D21").Value = 0.1105539106380888
$ws.Range("E21").Value = 0.03532101839856949
$ws.Range("G21").Value = 0.002581413973158118
$ws.Range("I21").Value = 2.480170681912384
$ws.Range("K21").Value = 1.801101829302922
$ws.Range("L21").Value = 0.3250327805373274
$ws.Range("N21").Value = 3.850922681232873
$ws.Range("B22").Value = 1.876058718375077
$ws.Range("C22").Value = 0.2032075803438147
$ws.Range("D22").Value = 0.1170700613010354
$ws.Range("E22").Value = 0.03555663196768855
$ws.Range("G22").Value = 0.002577780675061092
$ws.Range("I22").Value = 2.50336650112294
$ws.Range("K22").Value = 1.860137868114862
$ws.Range("L22").Value = 0.3332304670214228
$ws.Range("N22").Value = 3.866613819752445
$ws.Range("B23").Value = 1.84908891176093
$ws.Range("C23").Value = 0.1975549477419065
$ws.Range("D23").Value = 0.113588381330942
$ws.Range("E23").Value = 0.03542965990744662
$ws.Range("G23").Value = 0.002579707180022289
$ws.Range("I23").Value = 2.490918157694296
$ws.Range("K23").Value = 1.828547072354127
$ws.Range("L23").Value = 0.3288411432033485
$ws.Range("N23").Value = 3.858179856782499
$ws.Range("B24").Value = 1.748781405400223
$ws.Range("C24").Value = 0.176304834396916
$ws.Range("D24").Value = 0.1004866516970679
$ws.Range("E24").Value = 0.03497560785148224
$ws.Range("G24").Value = 0.002587282249342588
$ws.Range("I24").Value = 2.445269335446469
$ws.Range("K24").Value = 1.710700467596894
$ws.Range("L24").Value = 0.3125254351216995
$ws.Range("N24").Value = 3.827536945198858
$ws.Range("B25").Value = 1.644479860987872
$ws.Range("C25").Value = 0.1537106092235661
$ws.Range("D25").Value = 0.08652978800058975
$ws.Range("E25").Value = 0.03454328535397799
$ws.Range("G25").Value = 0.002596053363753071
$ws.Range("I25").Value = 2.399208891414673
$ws.Range("K25").Value = 1.587387057201795
$ws.Range("L25").Value = 0.2955795338682776
$ws.Range("N25").Value = 3.797237404604459
